# Rename the worksheet from "Property1" to "DataNode" (part of unifying the
# DataNode / DataTable / Entity concept across the config workbooks) and
# move the active selection to D26, matching the author's resulting state.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Name = "DataNode"

$ws.Activate() | Out-Null
$ws.Range("D26").Select() | Out-Null
